$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 43, pushing existing rows 43-129 down to 44-130.
$ws.Rows.Item(43).EntireRow.Insert()

# Populate the newly inserted row 43 with the new record's data.
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 45246
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 100112022
$ws.Range("G43").Value = "Arveja Verde"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 395
$ws.Range("K43").Value = 27000
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = 28405
$ws.Range("N43").Value = "$/saco 25 kilos"
$ws.Range("O43").Value = "Región del Maule"
$ws.Range("P43").Value = 1136
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"
